$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.222.70'
$ws.Range("E2").Value = '  +1.99%  '

$ws.Range("D3").Value = '3.120.75'
$ws.Range("E3").Value = '  +0.51%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = "'" + '237.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.30%  '

$ws.Range("D6").Value = "'" + '613.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D7").Value = "'" + '1.10'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.96%  '

$ws.Range("D8").Value = "'" + '0.391'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.79%  '

$ws.Range("D9").Value = "'" + '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Value = "'" + '0.841'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +13.24%  '

$ws.Range("D11").Value = '3.121.39'
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("E12").Value = '  -2.85%  '

$ws.Range("D13").Value = "'" + '0.0000245'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.34%  '

$ws.Range("D14").Value = "'" + '35.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.31%  '

$ws.Range("D15").Value = '93.073.81'
$ws.Range("E15").Value = '  +1.76%  '

$ws.Range("D16").Value = "'" + '5.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.00%  '

$ws.Range("D17").Value = '3.705.97'
$ws.Range("E17").Value = '  +0.62%  '

$ws.Range("D18").Value = '3.111.25'
$ws.Range("E18").Value = '  +0.68%  '

$ws.Range("D19").Value = "'" + '3.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("D20").Value = "'" + '14.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").Value = "'" + '5.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.20%  '

$ws.Range("E22").Value = '  +0.69%  '

$ws.Range("D23").Value = "'" + '442.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("D24").Value = "'" + '9.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.16%  '

$ws.Range("E25").Value = '  -1.41%  '

$ws.Range("D26").Value = "'" + '12.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.43%  '

$ws.Range("D27").Value = "'" + '86.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.59%  '

$ws.Range("D28").Value = '3.293.28'
$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +9.54%  '

$ws.Range("D31").Value = "'" + '0.238'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.13%  '

$ws.Range("D32").Value = "'" + '0.124'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -12.75%  '

$ws.Range("E33").Value = '  +4.42%  '

$ws.Range("D34").Value = "'" + '9.22'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.75%  '

$ws.Range("D35").Value = "'" + '8.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.13%  '

$ws.Range("D36").Value = "'" + '0.159'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.03%  '

$ws.Range("D37").Value = "'" + '25.99'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("D38").Value = "'" + '3.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.53%  '

$ws.Range("E39").Value = '  -1.50%  '

$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("D41").Value = "'" + '24.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.11%  '

$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").Value = "'" + '0.443'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.91%  '

$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = "'" + '474.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '

$ws.Range("E44").Value = '  -2.42%  '

$ws.Range("D46").Value = "'" + '158.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.48%  '

$ws.Range("D47").Value = "'" + '0.694'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("D49").Value = "'" + '1.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("E50").Value = '  +1.59%  '

$ws.Range("D51").Value = "'" + '44.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.15%  '

